$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q4 and shift the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# A7 is a brand-new cell (row 7 didn't exist before); copy the bold/border
# style used by the rest of column A (e.g. A6) onto it before writing the
# value.
$summary.Range("A6").Copy()
$summary.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2021-Q1"
$summary.Range("C7").Value = 2
$summary.Range("D7").Value = 0.42

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2021-Q2"
$summary.Range("C6").Value = 7
$summary.Range("D6").Value = 2.2

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 5
$summary.Range("D5").Value = 1.77

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 8
$summary.Range("D4").Value = 1.31

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 25
$summary.Range("D3").Value = 3.54

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 12
$summary.Range("D2").Value = 0.91

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet right before "2022-Q3", built from a
#    copy of "2022-Q3" so that all sheet-level formatting (header style,
#    outline props, column/row layout) matches the other quarterly sheets.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The copied sheet has 25 data rows (rows 2-26); the new quarter only has 12
# data rows (rows 2-13), so drop the surplus rows 14-26.
$q4.Range("A14:A26").EntireRow.Delete()

# Columns B-G hold text values (fund code / name / size / position figures
# stored as strings in the source data); force text formatting before
# writing so numeric-looking strings such as fund codes or "0.00" are not
# auto-converted to numbers.
$q4.Range("B2:G13").NumberFormat = "@"

$q4.Range("B2").Value = "164403"
$q4.Range("C2").Value = "前海开源沪港深农业混合（LOF）A"
$q4.Range("D2").Value = "4.16"
$q4.Range("E2").Value = "88.37"
$q4.Range("F2").Value = "5.10"
$q4.Range("G2").Value = "0.2122"
$q4.Range("H2").Value = 4

$q4.Range("B3").Value = "010418"
$q4.Range("C3").Value = "财通景气行业混合A"
$q4.Range("D3").Value = "2.60"
$q4.Range("E3").Value = "86.37"
$q4.Range("F3").Value = "7.48"
$q4.Range("G3").Value = "0.1945"
$q4.Range("H3").Value = 4

$q4.Range("B4").Value = "501015"
$q4.Range("C4").Value = "财通多策略升级混合（LOF）A"
$q4.Range("D4").Value = "1.97"
$q4.Range("E4").Value = "86.66"
$q4.Range("F4").Value = "7.38"
$q4.Range("G4").Value = "0.1454"
$q4.Range("H4").Value = 4

$q4.Range("B5").Value = "005959"
$q4.Range("C5").Value = "财通新视野灵活配置混合C"
$q4.Range("D5").Value = "1.61"
$q4.Range("E5").Value = "86.51"
$q4.Range("F5").Value = "7.61"
$q4.Range("G5").Value = "0.1225"
$q4.Range("H5").Value = 4

$q4.Range("B6").Value = "015210"
$q4.Range("C6").Value = "前海开源沪港深农业混合（LOF）C"
$q4.Range("D6").Value = "1.94"
$q4.Range("E6").Value = "88.37"
$q4.Range("F6").Value = "5.10"
$q4.Range("G6").Value = "0.0989"
$q4.Range("H6").Value = 4

$q4.Range("B7").Value = "005851"
$q4.Range("C7").Value = "财通新视野灵活配置混合A"
$q4.Range("D7").Value = "0.63"
$q4.Range("E7").Value = "86.51"
$q4.Range("F7").Value = "7.61"
$q4.Range("G7").Value = "0.0479"
$q4.Range("H7").Value = 4

$q4.Range("B8").Value = "015271"
$q4.Range("C8").Value = "财通多策略升级混合（LOF）C"
$q4.Range("D8").Value = "0.57"
$q4.Range("E8").Value = "86.66"
$q4.Range("F8").Value = "7.38"
$q4.Range("G8").Value = "0.0421"
$q4.Range("H8").Value = 4

$q4.Range("B9").Value = "010637"
$q4.Range("C9").Value = "财通安盈混合C"
$q4.Range("D9").Value = "1.58"
$q4.Range("E9").Value = "36.44"
$q4.Range("F9").Value = "1.60"
$q4.Range("G9").Value = "0.0253"
$q4.Range("H9").Value = 8

$q4.Range("B10").Value = "010636"
$q4.Range("C10").Value = "财通安盈混合A"
$q4.Range("D10").Value = "0.88"
$q4.Range("E10").Value = "36.44"
$q4.Range("F10").Value = "1.60"
$q4.Range("G10").Value = "0.0141"
$q4.Range("H10").Value = 8

$q4.Range("B11").Value = "006433"
$q4.Range("C11").Value = "平安鑫利灵活配置混合C"
$q4.Range("D11").Value = "0.23"
$q4.Range("E11").Value = "26.90"
$q4.Range("F11").Value = "0.95"
$q4.Range("G11").Value = "0.0022"
$q4.Range("H11").Value = 9

$q4.Range("B12").Value = "003626"
$q4.Range("C12").Value = "平安鑫利灵活配置混合A"
$q4.Range("D12").Value = "0.15"
$q4.Range("E12").Value = "26.90"
$q4.Range("F12").Value = "0.95"
$q4.Range("G12").Value = "0.0014"
$q4.Range("H12").Value = 9

$q4.Range("B13").Value = "016234"
$q4.Range("C13").Value = "财通景气行业混合C"
$q4.Range("D13").Value = "0.00"
$q4.Range("E13").Value = "86.37"
$q4.Range("F13").Value = "7.48"
$q4.Range("H13").Value = 4

# G13's source value is the literal number 0 (not text), unlike the other
# rows in this column, so clear the forced text format just for that cell
# before writing it as a real number.
$q4.Range("G13").NumberFormat = "General"
$q4.Range("G13").Value = 0
